# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (quarterly fund-holdings table) right
# after the "总计" (totals) summary sheet, fills it with the Q3 data, and
# updates the "总计" summary table with a new leading row for 2022-Q3
# (pushing the existing quarterly rows down by one).

$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell while forcing it to be stored as TEXT
# (never auto-converted to a number by the "smart" Value setter), and
# without leaving behind any NumberFormat/style residue on the target
# cell. We stage the value in a scratch cell far outside the used range,
# force that scratch cell to Text format, assign the value there, then
# copy only the *value* (PasteSpecial values) across to the real target.
function Set-TextValue($ws, $addr, $val) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $val
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $scratch.Clear()
}

$totalSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet.
#    We clone the existing "2022-Q2" sheet (currently Worksheets.Item(2))
#    so the new sheet inherits identical structure/styles, place the
#    clone right after "总计", rename it, drop the extra data row that
#    2022-Q3 doesn't need, and overwrite its contents with the Q3 data.
# ---------------------------------------------------------------------

$templateSheet = $wb.Worksheets.Item(2)
$templateSheet.Copy($null, $totalSheet)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Template had 3 data rows (rows 2-4); 2022-Q3 only needs 2 (rows 2-3).
$q3Sheet.Rows.Item(4).Delete()

# Row 2: 164811 / 工银瑞信中证京津冀协同发展主题指数（LOF）A
Set-TextValue $q3Sheet "B2" "164811"
$q3Sheet.Range("C2").Value = "工银瑞信中证京津冀协同发展主题指数（LOF）A"
Set-TextValue $q3Sheet "D2" "0.12"
Set-TextValue $q3Sheet "E2" "93.09"
Set-TextValue $q3Sheet "F2" "3.13"
Set-TextValue $q3Sheet "G2" "0.0038"
$q3Sheet.Range("H2").Value = 3

# Row 3: 164825 / 工银瑞信中证京津冀协同发展主题指数（LOF）C
Set-TextValue $q3Sheet "B3" "164825"
$q3Sheet.Range("C3").Value = "工银瑞信中证京津冀协同发展主题指数（LOF）C"
Set-TextValue $q3Sheet "D3" "0.03"
Set-TextValue $q3Sheet "E3" "93.09"
Set-TextValue $q3Sheet "F3" "3.13"
Set-TextValue $q3Sheet "G3" "0.0009"
$q3Sheet.Range("H3").Value = 3

# ---------------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert a new leading data row for
#    2022-Q3 and push the previously existing quarters down by one row.
# ---------------------------------------------------------------------

$totalSheet.Range("A7").Value = 5
$totalSheet.Range("B7").Value = "2020-Q4"
$totalSheet.Range("C7").Value = 3
$totalSheet.Range("D7").Value = 0.02
# A7 is a brand-new cell beyond the old A1:D6 range; give it the same
# style as the other index cells in column A (bold/centered/bordered).
$totalSheet.Range("A6").Copy()
$totalSheet.Range("A7").PasteSpecial(-4122)
$totalSheet.Range("A7").Value = 5

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q3"
$totalSheet.Range("C6").Value = 3
$totalSheet.Range("D6").Value = 0.01

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q4"
$totalSheet.Range("C5").Value = 3
$totalSheet.Range("D5").Value = 0.01

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q1"
$totalSheet.Range("C4").Value = 4
$totalSheet.Range("D4").Value = 0.02

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 3
$totalSheet.Range("D3").Value = 0.01

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0

# ---------------------------------------------------------------------
# 3) Restore the originally-active tab ("2020-Q4", the last sheet) since
#    copying a sheet shifts Excel's active-sheet selection to the copy.
# ---------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
